$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT (keeps leading
# zeroes / decimal-look-alike strings from being coerced to numbers),
# while leaving the cell on the default (unstyled) format - mirrors the
# source data which stores these as inlineStr with no explicit style.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (position 2)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q4Sheet)
$newSheet.Name = "2022-Q3"

# Header row (basic text values first)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 (fund: 014831 / 兴银中证1000指数增强A)
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "014831"
$newSheet.Range("C2").Value = "兴银中证1000指数增强A"
Set-TextValue $newSheet.Range("D2") "1.37"
Set-TextValue $newSheet.Range("E2") "83.33"
Set-TextValue $newSheet.Range("F2") "1.08"
Set-TextValue $newSheet.Range("G2") "0.0148"
$newSheet.Range("H2").Value = 4

# Row 3 (fund: 014832 / 兴银中证1000指数增强C)
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "014832"
$newSheet.Range("C3").Value = "兴银中证1000指数增强C"
Set-TextValue $newSheet.Range("D3") "0.90"
Set-TextValue $newSheet.Range("E3") "83.33"
Set-TextValue $newSheet.Range("F3") "1.08"
Set-TextValue $newSheet.Range("G3") "0.0097"
$newSheet.Range("H3").Value = 4

# Re-apply number formats back to General so the only thing that stuck
# from the "@" trick is the text storage type, not a lingering style.
$newSheet.Range("Z100").Copy()
$newSheet.Range("D2:G3").PasteSpecial(-4122)
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Range("B3").PasteSpecial(-4122)

# Apply the bold/bordered header style (copied from the "总计" header)
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert a new summary row into "总计" (row 2) for 2022-Q3, pushing
#    the existing 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 1.62

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 1.68

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.02

# ---------------------------------------------------------------------
# 3. Restore the original active-sheet bookkeeping: before this edit,
#    "2021-Q3" was the selected/active tab - keep it that way instead of
#    leaving the freshly-inserted sheet active.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$q3Sheet.Activate()
